$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the "Modified Reg iProctor P2,P3" data
$ws.Range("A2").Value = "xdCzd476"
$ws.Range("B2").Value = 23081409
$ws.Range("C2").Value = "dflgtzp77"
$ws.Range("D2").Value = "j#8SkF`$6"
$ws.Range("F2").Value = "CKywQwIb"
$ws.Range("G2").Value = "IKYx"
